$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.142.36"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.15%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.669.64"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.58%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.39%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5224"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.59%  "
$ws.Range("E7").Value = "  -0.39%  "
$ws.Range("E8").Value = "  -2.22%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06328"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.43%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.21"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.86%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07551"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.84%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.682.29"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.28%  "
$ws.Range("E13").Value = "  -1.98%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5487"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008031"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.48"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.23%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.173.58"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("E18").Value = "  -0.32%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.755"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "187.78"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.31"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.76%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.241"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.38%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.003"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.35%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "149.36"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.20%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1243"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.94%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.482"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.81%  "
$ws.Range("E27").Value = "  -0.40%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06315"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.59%  "
$ws.Range("E29").Value = "  -1.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.284"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.530"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.419"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.650"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.46%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.006"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.81%  "
$ws.Range("E35").Value = "  -0.98%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.764"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.72%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.394"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.119.34"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.46%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01612"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.072"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.44%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8647"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.80%  "
$ws.Range("E42").Value = "  -0.66%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.40"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.822.47"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.58%  "
$ws.Range("E45").Value = "  -0.85%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "55.55"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.22%  "
$ws.Range("E47").Value = "  +0.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.068"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05234"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.58%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4237"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.925"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.17%  "
